$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire row 12 (as if the user clicked the row header) and
# clear its contents (Delete key / Clear Contents), same as the row 12
# data disappearing from the sheet in the target file.
$ws.Range("A12:XFD12").Select()
$ws.Rows.Item(12).ClearContents()
